$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column J header (row 5) - text with a line break, like the other headers
$ws.Range("J5").Value = "No standardization" + [char]10 + "shuffle when use kfold (seed = 1)"

# Copy cell formatting (number format / alignment) from sibling cells so the
# new column reuses the existing style entries instead of creating new ones
$ws.Range("H5").Copy() | Out-Null
$ws.Range("J5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B6").Copy() | Out-Null
$ws.Range("J6:J10").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# New column J data values (rows 6:10)
$ws.Range("J6").Value = 4.8006113476114498
$ws.Range("J7").Value = 4.8002483703501797
$ws.Range("J8").Value = 4.8127441256846701
$ws.Range("J9").Value = 4.8634390971983104
$ws.Range("J10").Value = 4.9566232289732701

# Widen the new column to fit the header text
$ws.Columns.Item(10).ColumnWidth = 35 + 2/3

# Move the view / selection over to the new column, matching the saved view
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("J18").Select() | Out-Null
